# hautes herbes avec differents pokemons
#
# Adds a "Categorie" column to both the "element_terrain" and "pokemon"
# sheets, and adds a new "HAUTE_HERBE" terrain row + a new "Dracofeu"
# pokemon row.

$wb  = $excel.ActiveWorkbook
$ws3 = $wb.Worksheets.Item("element_terrain")
$ws4 = $wb.Worksheets.Item("pokemon")

# ---------------------------------------------------------------------
# element_terrain: insert a new "Categorie" column at column C, shifting
# the existing "traversable" / "% proba" columns one to the right
# (C->D, D->E). Columns are shifted right-to-left via copy/paste so the
# original column width metadata for C/D is left untouched (matches the
# target file, where the new E column simply gets its own <col> entry).
# ---------------------------------------------------------------------
for ($c = 4; $c -ge 3; $c--) {
    $src = $ws3.Range($ws3.Cells.Item(1, $c), $ws3.Cells.Item(4, $c))
    $dstCol = $c + 1
    $dst = $ws3.Range($ws3.Cells.Item(1, $dstCol), $ws3.Cells.Item(4, $dstCol))
    $src.Copy() | Out-Null
    $dst.PasteSpecial(-4104) | Out-Null
}

$ws3.Range("C1").Value = "Categorie"
$ws3.Range("C2").Value = "DEBUT"
$ws3.Range("C3").Value = "DEBUT"
$ws3.Range("C4").Value = "DEBUT"

$ws3.Columns.Item(5).ColumnWidth = 38.42578125

# New row 5: a "HAUTE_HERBE" terrain entry
$ws3.Range("A5").Value = "HAUTE_HERBE"
$ws3.Range("B5").Value = "path_herbe"
$ws3.Range("C5").Value = "HARD"
$ws3.Range("D5").Value = "OUI"
$ws3.Range("E5").Value = 50

# A5 gets a distinct look: Consolas font, vertically centred
$ws3.Range("A5").Font.Name = "Consolas"
$ws3.Range("A5").VerticalAlignment = -4108

# A blank placeholder cell also appears in column C of row 8
$ws3.Range("A8").Copy() | Out-Null
$ws3.Range("C8").PasteSpecial(-4122) | Out-Null

$ws3.PageSetup.PaperSize = 9
$ws3.PageSetup.Orientation = 1

# ---------------------------------------------------------------------
# pokemon: insert a new "Categorie" column at column C, shifting the
# existing "rarete/atk/def/pvMax" columns one to the right (C->D, D->E,
# E->F, F->G).
# ---------------------------------------------------------------------
for ($c = 6; $c -ge 3; $c--) {
    $src = $ws4.Range($ws4.Cells.Item(1, $c), $ws4.Cells.Item(3, $c))
    $dstCol = $c + 1
    $dst = $ws4.Range($ws4.Cells.Item(1, $dstCol), $ws4.Cells.Item(3, $dstCol))
    $src.Copy() | Out-Null
    $dst.PasteSpecial(-4104) | Out-Null
}

$ws4.Range("C1").Value = "Categorie"
$ws4.Range("C2").Value = "DEBUT"
$ws4.Range("C3").Value = "DEBUT"

# New row 4: a "Dracofeu" pokemon entry
$ws4.Range("A4").Value = "Dracofeu"
$ws4.Range("B4").Value = "path_dracofeu"
$ws4.Range("C4").Value = "HARD"
$ws4.Range("D4").Value = "100.0"
$ws4.Range("E4").Value = 200
$ws4.Range("F4").Value = 200
$ws4.Range("G4").Value = 2000

# ---------------------------------------------------------------------
# Selections / active sheet: "pokemon" stays the active tab, with D4
# selected; "element_terrain" keeps a selection on E5.
# ---------------------------------------------------------------------
$ws3.Activate()
$ws3.Range("E5").Select()

$ws4.Activate()
$ws4.Range("D4").Select()
